$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows are rotated: row2's original values move to row4,
# row3's original values move to row2, row4's original values move to row3.
# Capture the original values first, then write the rotated values.

$cols = @("D", "J", "K", "L", "M", "P")

$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("${col}2").Value2
    $orig3[$col] = $ws.Range("${col}3").Value2
    $orig4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $orig4[$col]
    $ws.Range("${col}3").Value2 = $orig2[$col]
    $ws.Range("${col}4").Value2 = $orig3[$col]
}
